$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 11, pushing the existing
# rows 11-20 down to 12-21 (dimension grows from A1:T20 to A1:T21).
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C11").Value = 'Arica y Parinacota'
$ws.Range("D11").Value = 44650
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 'Fruta'
$ws.Range("G11").Value = 100103
$ws.Range("H11").Value = 'Frutos de hueso (carozo)'
$ws.Range("I11").Value = 100103002
$ws.Range("J11").Value = 'Ciruela'
$ws.Range("K11").Value = 'Angeleno'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 17000
$ws.Range("O11").Value = 18000
$ws.Range("P11").Value = 17500
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 972
$ws.Range("T11").Value = 18
